$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style/format of the row above (row 16) into the new row 17
# so the date cell (G17) picks up the existing date style (s="1")
# instead of minting a brand new number format.
$ws.Range("A16:H16").Copy()
$ws.Range("A17:H17").PasteSpecial(-4122)

$ws.Range("A17").Value = 9300.5300000000007
$ws.Range("B17").Value = 9578.2999999999993
$ws.Range("C17").Value = 294.14
$ws.Range("D17").Value = 302.66000000000003
$ws.Range("E17").Value = $true
$ws.Range("F17").Value = 2.9
$ws.Range("G17").Value = 42626.544317129628
$ws.Range("H17").Value = $false
